$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Header row for the vehicle ("载具") data table
$ws.Range("A1").Value = "载具编号"
$ws.Range("B1").Value = "载具名称"
$ws.Range("C1").Value = "载具类型"
$ws.Range("D1").Value = "额外技能"
$ws.Range("E1").Value = "获取途径"
$ws.Range("F1").Value = "评分"

# Leave the selection where it would naturally land after typing the last header
$ws.Range("G1").Select()
